# Weekly update: insert a new Choclo price record for
# "Comercializadora del Agro de Limarí" (Coquimbo), pushing the
# existing rows 201-203 down to 202-204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 201; this shifts the
# old rows 201, 202, 203 down to 202, 203, 204 respectively (carrying
# their formatting, e.g. the date style on column D, along with them).
$ws.Range("A201").EntireRow.Insert()

# Populate the newly inserted row 201 with the new weekly record.
$ws.Range("A201").Value = 2
$ws.Range("B201").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C201").Value = "Coquimbo"
$ws.Range("D201").Value = 45239
$ws.Range("E201").Value = 4
$ws.Range("F201").Value = 100112024
$ws.Range("G201").Value = "Choclo"
$ws.Range("H201").Value = "Dulce o Americano"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 600
$ws.Range("K201").Value = 26000
$ws.Range("L201").Value = 28000
$ws.Range("M201").Value = 27000
$ws.Range("N201").Value = "`$/malla 70 unidades"
$ws.Range("O201").Value = "Provincia de Limarí"
$ws.Range("P201").Value = 386
$ws.Range("Q201").Value = 70
$ws.Range("R201").Value = "Hortaliza"
